# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Column D mixes "thousand dot" formatted numbers (e.g. "25.762.76") with
# plain decimals (e.g. "0.0632"); the workbook stores the whole column as
# text, so values that parse as a plain number are entered with a leading
# apostrophe to force Excel to keep them as text instead of coercing them
# into a Number cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "25.762.76"
$ws.Cells.Item(2, 5).Value = "  -0.11%  "
$ws.Cells.Item(3, 4).Value = "1.631.12"
$ws.Cells.Item(3, 5).Value = "  -0.24%  "
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
$ws.Cells.Item(5, 4).Value = "'215.12"
$ws.Cells.Item(5, 5).Value = "  -0.10%  "
$ws.Cells.Item(6, 5).Value = "  -0.64%  "
$ws.Cells.Item(7, 5).Value = "  -0.09%  "
$ws.Cells.Item(8, 5).Value = "  -0.63%  "
$ws.Cells.Item(9, 4).Value = "'0.0632"
$ws.Cells.Item(9, 5).Value = "  -1.45%  "
$ws.Cells.Item(10, 4).Value = "'19.47"
$ws.Cells.Item(10, 5).Value = "  -1.83%  "
$ws.Cells.Item(11, 4).Value = "'0.0793"
$ws.Cells.Item(11, 5).Value = "  +0.85%  "
$ws.Cells.Item(12, 5).Value = "  +0.27%  "
$ws.Cells.Item(13, 4).Value = "1.857.34"
$ws.Cells.Item(13, 5).Value = "  -0.14%  "
$ws.Cells.Item(14, 4).Value = "1.631.59"
$ws.Cells.Item(14, 5).Value = "  -0.55%  "
$ws.Cells.Item(15, 4).Value = "'0.556"
$ws.Cells.Item(15, 5).Value = "  +0.18%  "
$ws.Cells.Item(16, 5).Value = "  -1.67%  "
$ws.Cells.Item(17, 4).Value = "'63.07"
$ws.Cells.Item(17, 5).Value = "  -0.03%  "
$ws.Cells.Item(18, 4).Value = "25.769.00"
$ws.Cells.Item(18, 5).Value = "  -0.13%  "
$ws.Cells.Item(19, 5).Value = "  -0.11%  "
$ws.Cells.Item(20, 5).Value = "  -0.08%  "
$ws.Cells.Item(21, 4).Value = "'192.29"
$ws.Cells.Item(21, 5).Value = "  -1.10%  "
$ws.Cells.Item(22, 5).Value = "  -0.08%  "
$ws.Cells.Item(23, 5).Value = "  +1.88%  "
$ws.Cells.Item(24, 5).Value = "  -0.10%  "
$ws.Cells.Item(25, 5).Value = "  +2.91%  "
$ws.Cells.Item(26, 4).Value = "'142.81"
$ws.Cells.Item(26, 5).Value = "  +2.22%  "
$ws.Cells.Item(27, 5).Value = "  +1.97%  "
$ws.Cells.Item(28, 4).Value = "'6.86"
$ws.Cells.Item(28, 5).Value = "  +0.57%  "
$ws.Cells.Item(29, 5).Value = "  -0.58%  "
$ws.Cells.Item(30, 5).Value = "  -0.14%  "
$ws.Cells.Item(31, 5).Value = "  -0.59%  "
$ws.Cells.Item(33, 4).Value = "'3.22"
$ws.Cells.Item(33, 5).Value = "  -0.71%  "
$ws.Cells.Item(34, 5).Value = "  -1.37%  "
$ws.Cells.Item(35, 5).Value = "  -0.46%  "
$ws.Cells.Item(36, 4).Value = "'0.901"
$ws.Cells.Item(36, 5).Value = "  +0.46%  "
$ws.Cells.Item(37, 4).Value = "1.133.40"
$ws.Cells.Item(37, 5).Value = "  +2.01%  "
$ws.Cells.Item(38, 5).Value = "  -2.11%  "
$ws.Cells.Item(39, 5).Value = "  -1.43%  "
$ws.Cells.Item(40, 5).Value = "  -1.12%  "
$ws.Cells.Item(41, 5).Value = "  +0.07%  "
$ws.Cells.Item(42, 5).Value = "  +0.71%  "
$ws.Cells.Item(43, 2).Value = "Quant"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(43, 4).Value = "'100.60"
$ws.Cells.Item(43, 5).Value = "  +1.34%  "
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).Value = "'5.53"
$ws.Cells.Item(44, 5).Value = "  -0.79%  "
$ws.Cells.Item(45, 5).Value = "  -0.64%  "
$ws.Cells.Item(46, 4).Value = "1.766.91"
$ws.Cells.Item(46, 5).Value = "  -0.06%  "
$ws.Cells.Item(47, 2).Value = "Aave"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(47, 4).Value = "'55.30"
$ws.Cells.Item(47, 5).Value = "  -0.46%  "
$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48, 4).Value = "'0.0507"
$ws.Cells.Item(48, 5).Value = "  +0.76%  "
$ws.Cells.Item(49, 2).Value = "Mantle"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(49, 4).Value = "'0.418"
$ws.Cells.Item(49, 5).Value = "  +0.09%  "
$ws.Cells.Item(50, 2).Value = "SynthetixNetwork"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Cells.Item(50, 4).Value = "'2.35"
$ws.Cells.Item(50, 5).Value = "  -6.81%  "
$ws.Cells.Item(51, 2).Value = "RenderToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(51, 4).Value = "'1.41"
$ws.Cells.Item(51, 5).Value = "  +2.75%  "
